# Normalize the "Recorded By" (column G) values on the active sheet:
# move the "System"/"system" token from the front of the comma-separated
# list to the end, for the specific values that need it.
#
# Mapping of exact cell text -> new exact cell text
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$replacements = @{
    "System, system, backup@backdoor.com" = "System, backup@backdoor.com, system"
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "System, admin@admin.com"             = "admin@admin.com, System"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2
    if ($null -ne $val -and $replacements.ContainsKey($val)) {
        $cell.Value = $replacements[$val]
    }
}
